# TC21_Verify_ Find_a_Branch_Guestuser.xlsx
# "KFP Support files added & updated"
#
# The test case's "Search by Zip/State" flow had a redundant verification
# step removed (the "ZipcodeBranches" web-element check right after the
# zip-code search), and the state used for the "search by state" branch of
# the test was swapped from Maryland/20152 to Idaho/83815, with the
# "Find a Branch Near You" banner text corrected to "Find a Location Near
# You" on the Testdata sheet.

$wb = $excel.ActiveWorkbook

$wsTest  = $wb.Worksheets.Item(1)   # TC21_Verify_Find_a_Branch_Guser
$wsData  = $wb.Worksheets.Item(2)   # Testdata

# Remove the two rows for the now-dropped "ZipcodeBranches" verification
# step (VERIFY_WEBELEMENT_PRESENT + its trailing WAIT). Everything below
# shifts up by two rows automatically.
$wsTest.Rows("13:14").Delete()

# Update the Testdata key/value pairs used by the steps above.
$wsData.Range("B4").Value = "Find a Location Near You"
$wsData.Range("B6").Value = 83815
$wsData.Range("B11").Value = "Idaho"

# Restore view selections to match the authored state: Testdata's selection
# first (it is not the active sheet), then re-activate the test-case sheet
# and select its final range so it remains the active tab.
[void]$wsData.Range("B12").Select()
$wsTest.Activate()
[void]$wsTest.Range("A13:XFD14").Select()
